$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 20.016918182373047
$ws.Range("C2").Value = 5.712643623352051
$ws.Range("D2").Value = 13.816729545593262
$ws.Range("E2").Value = 57.85714340209961
